{"js": "// Replace each two-digit multiplication expression with its new value.\n// Every \"before\" string below is unique in the document, so a plain\n// text search + Replace is unambiguous.\nconst replacements = [\n  [\"71\u00d799=\", \"15\u00d751=\"],\n  [\"63\u00d746=\", \"50\u00d795=\"],\n  [\"42\u00d721=\", \"76\u00d717=\"],\n  [\"72\u00d751=\", \"78\u00d727=\"],\n  [\"27\u00d714=\", \"32\u00d713=\"],\n  [\"61\u00d758=\", \"76\u00d715=\"],\n  [\"15\u00d729=\", \"42\u00d736=\"],\n  [\"12\u00d730=\", \"20\u00d780=\"],\n  [\"12\u00d768=\", \"65\u00d766=\"],\n  [\"84\u00d753=\", \"43\u00d769=\"],\n  [\"94\u00d766=\", \"33\u00d791=\"],\n  [\"98\u00d795=\", \"56\u00d757=\"],\n  [\"84\u00d711=\", \"28\u00d780=\"],\n  [\"78\u00d793=\", \"49\u00d763=\"],\n  [\"95\u00d793=\", \"51\u00d720=\"],\n  [\"67\u00d784=\", \"81\u00d769=\"],\n  [\"50\u00d734=\", \"82\u00d756=\"],\n  [\"85\u00d745=\", \"18\u00d760=\"],\n  [\"82\u00d719=\", \"78\u00d796=\"],\n  [\"77\u00d720=\", \"64\u00d748=\"],\n  [\"11\u00d745=\", \"81\u00d759=\"],\n  [\"66\u00d798=\", \"22\u00d760=\"],\n  [\"71\u00d740=\", \"72\u00d744=\"],\n  [\"32\u00d743=\", \"84\u00d724=\"],\n  [\"20\u00d782=\", \"33\u00d794=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const found = body.search(before, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication expression with its new value.\n# Every \"before\" string is unique in the document, so Find/Replace with\n# wdReplaceAll for each pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"71\u00d799=\", \"15\u00d751=\"),\n    @(\"63\u00d746=\", \"50\u00d795=\"),\n    @(\"42\u00d721=\", \"76\u00d717=\"),\n    @(\"72\u00d751=\", \"78\u00d727=\"),\n    @(\"27\u00d714=\", \"32\u00d713=\"),\n    @(\"61\u00d758=\", \"76\u00d715=\"),\n    @(\"15\u00d729=\", \"42\u00d736=\"),\n    @(\"12\u00d730=\", \"20\u00d780=\"),\n    @(\"12\u00d768=\", \"65\u00d766=\"),\n    @(\"84\u00d753=\", \"43\u00d769=\"),\n    @(\"94\u00d766=\", \"33\u00d791=\"),\n    @(\"98\u00d795=\", \"56\u00d757=\"),\n    @(\"84\u00d711=\", \"28\u00d780=\"),\n    @(\"78\u00d793=\", \"49\u00d763=\"),\n    @(\"95\u00d793=\", \"51\u00d720=\"),\n    @(\"67\u00d784=\", \"81\u00d769=\"),\n    @(\"50\u00d734=\", \"82\u00d756=\"),\n    @(\"85\u00d745=\", \"18\u00d760=\"),\n    @(\"82\u00d719=\", \"78\u00d796=\"),\n    @(\"77\u00d720=\", \"64\u00d748=\"),\n    @(\"11\u00d745=\", \"81\u00d759=\"),\n    @(\"66\u00d798=\", \"22\u00d760=\"),\n    @(\"71\u00d740=\", \"72\u00d744=\"),\n    @(\"32\u00d743=\", \"84\u00d724=\"),\n    @(\"20\u00d782=\", \"33\u00d794=\")\n)\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
